$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.095.10'
$ws.Range("E2").Value = '  +1.86%  '

$ws.Range("D3").Value = '3.109.38'
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.84'
$ws.Range("E5").Value = '  -0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '615.03'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("E7").Value = '  -1.70%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.389'
$ws.Range("E8").Value = '  +4.84%  '

$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").Value = '3.108.48'
$ws.Range("E10").Value = '  -1.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.736'
$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000247'
$ws.Range("E13").Value = '  +0.66%  '

$ws.Range("D14").Value = '92.245.40'
$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.07'
$ws.Range("E15").Value = '  -2.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.43'
$ws.Range("E16").Value = '  -1.76%  '

$ws.Range("D17").Value = '3.698.26'
$ws.Range("E17").Value = '  -0.85%  '

$ws.Range("D18").Value = '3.112.16'
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.72'
$ws.Range("E19").Value = '  +0.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.62'
$ws.Range("E20").Value = '  -2.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.76'
$ws.Range("E21").Value = '  -3.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.37'
$ws.Range("E22").Value = '  +3.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '446.13'
$ws.Range("E23").Value = '  -0.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000195'
$ws.Range("E24").Value = '  -3.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.74'
$ws.Range("E25").Value = '  -3.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.58'
$ws.Range("E26").Value = '  -2.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.72'
$ws.Range("E27").Value = '  -1.17%  '

$ws.Range("D28").Value = '3.273.90'
$ws.Range("E28").Value = '  -1.63%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.132'
$ws.Range("E30").Value = '  -4.88%  '

$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.170'
$ws.Range("E32").Value = '  -0.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.10'
$ws.Range("E33").Value = '  -1.45%  '

$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.992'
$ws.Range("E34").Value = '  -1.02%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.86'
$ws.Range("E35").Value = '  +3.86%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.158'
$ws.Range("E36").Value = '  -8.36%  '

$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.13'
$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.90'
$ws.Range("E38").Value = '  -3.45%  '

$ws.Range("B39").Value = 'MantraDAO'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.88'
$ws.Range("E39").Value = '  +1.31%  '

$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '483.00'
$ws.Range("E40").Value = '  -4.36%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.29'
$ws.Range("E41").Value = '  -2.76%  '

$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.87'
$ws.Range("E42").Value = '  +8.16%  '

$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.432'
$ws.Range("E43").Value = '  -2.75%  '

$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.28'
$ws.Range("E44").Value = '  -4.00%  '

$ws.Range("E45").Value = '  +0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '161.99'
$ws.Range("E46").Value = '  +4.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.690'
$ws.Range("E47").Value = '  -3.22%  '

$ws.Range("E48").Value = '  -1.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.39'
$ws.Range("E49").Value = '  +1.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0335'
$ws.Range("E50").Value = '  +4.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.41'
$ws.Range("E51").Value = '  -0.46%  '

$ws.Range("D2:E51").ClearFormats()
